$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in this exact order so their
# --- sharedStrings.xml indices come out as 8 (space), 9 (Transferwise), 10 (BPCS) ---
$ws.Range("C9").Value = " "
$ws.Range("C8").Value = "Leí sobre Transferwise. Genial está. Probablemente de ahí podré transferir fondos y recibir mis ganancias de TradeStation. Y de Zulutrade en Skrill. Y después, todo a Payoneer. "
$ws.Range("C7").Value = "Vi que no me ejecutó cierta orden de BPCS de TSLA. Mary me dijo que podría ser por mi nivel 1 en TS. Tengo que ver para aumentar sus fondos. "

# --- Dates for the two new rows, re-using the existing date number format (style index 2) ---
$ws.Range("B5").Copy()
$ws.Range("B7:B8").PasteSpecial(-4122)
$ws.Range("B7").Value = 44080
$ws.Range("B8").Value = 44110

# --- Column B width (so it isn't left at the default) ---
$ws.Columns("B").ColumnWidth = 8.86

# --- Update selection to match the new state ---
$ws.Range("C14").Select()
